$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.874.18'
$ws.Range('E2').Value = '  +0.55%  '
$ws.Range('D3').Value = '2.237.45'
$ws.Range('E3').Value = '  +1.88%  '
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').Formula = "'270.10"
$ws.Range('E5').Value = '  +3.78%  '
$ws.Range('D6').Formula = "'94.16"
$ws.Range('E6').Value = '  +14.70%  '
$ws.Range('D7').Formula = "'0.629"
$ws.Range('E7').Value = '  +1.15%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').Formula = "'0.638"
$ws.Range('E9').Value = '  +7.79%  '
$ws.Range('D10').Formula = "'46.11"
$ws.Range('E10').Value = '  +5.76%  '
$ws.Range('E11').Value = '  +4.15%  '
$ws.Range('D12').Formula = "'8.41"
$ws.Range('E12').Value = '  +20.90%  '
$ws.Range('E13').Value = '  +1.88%  '
$ws.Range('E14').Value = '  +7.18%  '
$ws.Range('D15').Value = '2.573.67'
$ws.Range('E15').Value = '  +2.12%  '
$ws.Range('D16').Formula = "'0.819"
$ws.Range('E16').Value = '  +5.11%  '
$ws.Range('D17').Value = '2.228.13'
$ws.Range('E17').Value = '  +0.69%  '
$ws.Range('D18').Value = '43.887.35'
$ws.Range('E18').Value = '  +0.80%  '
$ws.Range('E19').Value = '  +1.71%  '
$ws.Range('E20').Value = '  +4.69%  '
$ws.Range('D21').Formula = "'70.82"
$ws.Range('E21').Value = '  +1.43%  '
$ws.Range('D22').Formula = "'2.34"
$ws.Range('E22').Value = '  -4.15%  '
$ws.Range('D23').Formula = "'234.67"
$ws.Range('E23').Value = '  +1.81%  '
$ws.Range('D24').Formula = "'9.10"
$ws.Range('E24').Value = '  +2.62%  '
$ws.Range('E25').Value = '  -0.04%  '
$ws.Range('D26').Formula = "'11.41"
$ws.Range('E26').Value = '  +6.82%  '
$ws.Range('E27').Value = '  +11.82%  '
$ws.Range('E28').Value = '  +6.26%  '
$ws.Range('D29').Formula = "'40.44"
$ws.Range('E29').Value = '  -4.33%  '
$ws.Range('E30').Value = '  +2.69%  '
$ws.Range('D31').Formula = "'172.73"
$ws.Range('E31').Value = '  -0.67%  '
$ws.Range('E32').Value = '  +5.25%  '
$ws.Range('D33').Formula = "'21.02"
$ws.Range('E33').Value = '  +2.83%  '
$ws.Range('E34').Value = '  +3.34%  '
$ws.Range('E35').Value = '  +2.15%  '
$ws.Range('E36').Value = '  -1.01%  '
$ws.Range('D37').Formula = "'0.0353"
$ws.Range('E37').Value = '  +0.18%  '
$ws.Range('E38').Value = '  -3.77%  '
$ws.Range('D39').Formula = "'3.57"
$ws.Range('E39').Value = '  +23.94%  '
$ws.Range('D40').Formula = "'12.80"
$ws.Range('E40').Value = '  -2.33%  '
$ws.Range('E41').Value = '  +12.46%  '
$ws.Range('D42').Formula = "'2.16"
$ws.Range('E42').Value = '  +3.13%  '
$ws.Range('D43').Formula = "'63.37"
$ws.Range('E43').Value = '  -0.75%  '
$ws.Range('D44').Formula = "'5.41"
$ws.Range('E44').Value = '  -0.88%  '
$ws.Range('D45').Formula = "'0.0997"
$ws.Range('E45').Value = '  +1.69%  '
$ws.Range('D46').Formula = "'101.35"
$ws.Range('E46').Value = '  +1.03%  '
$ws.Range('D47').Formula = "'8.40"
$ws.Range('E47').Value = '  +1.65%  '
$ws.Range('E48').Value = '  +4.26%  '
$ws.Range('E49').Value = '  +2.38%  '
$ws.Range('D50').Formula = "'0.449"
$ws.Range('E50').Value = '  +2.74%  '
$ws.Range('D51').Value = '2.458.51'
$ws.Range('E51').Value = '  +2.01%  '
